$d = $word.ActiveDocument

# Locate the paragraph that ends with the FLEMMING bibliography entry, then
# remove everything from the (empty) paragraph right after it through the
# end of the "© 2020 ..." paragraph: that is, the blank separator paragraph,
# the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph, and the
# "© 2020 . Contact: ..." paragraph. The blank paragraph that originally
# followed the copyright line is left intact, as is the page-break paragraph
# after it.

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*FLEMMING, Diva M.*Pearson Prentice Hall, 2009.*") {
        $target = $i
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the FLEMMING bibliography paragraph"
}

$startPara = $d.Paragraphs.Item($target + 1)
$endPara = $d.Paragraphs.Item($target + 3)

$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()
